# Auto-generated Excel COM-interop edit script
# Applies numeric "want-to-go" count bumps across sheets 1-4,
# plus a content re-shuffle (new row inserted + stale row dropped) in sheet 4 rows 36-49.

$wb = $excel.ActiveWorkbook

# Sheet 1 (展览): F-column (want-to-go count) updates
$ws = $wb.Worksheets.Item(1)
$ws.Range("F4").Value = 1144
$ws.Range("F6").Value = 154
$ws.Range("F10").Value = 1225
$ws.Range("F11").Value = 27671
$ws.Range("F12").Value = 3125
$ws.Range("F14").Value = 225
$ws.Range("F15").Value = 433
$ws.Range("F17").Value = 286
$ws.Range("F19").Value = 254
$ws.Range("F20").Value = 226
$ws.Range("F21").Value = 330
$ws.Range("F22").Value = 18
$ws.Range("F26").Value = 476
$ws.Range("F29").Value = 568
$ws.Range("F30").Value = 225
$ws.Range("F31").Value = 28

# Sheet 2 (演出): F-column (want-to-go count) updates
$ws = $wb.Worksheets.Item(2)
$ws.Range("F7").Value = 691
$ws.Range("F8").Value = 76
$ws.Range("F10").Value = 4209
$ws.Range("F15").Value = 35
$ws.Range("F21").Value = 4208

# Sheet 3 (本地生活): F-column (want-to-go count) updates
$ws = $wb.Worksheets.Item(3)
$ws.Range("F4").Value = 1098

# Sheet 4 (全部类型): F-column (want-to-go count) updates
$ws = $wb.Worksheets.Item(4)
$ws.Range("F4").Value = 1098
$ws.Range("F11").Value = 691
$ws.Range("F13").Value = 1144
$ws.Range("F14").Value = 154
$ws.Range("F18").Value = 1225
$ws.Range("F19").Value = 27672
$ws.Range("F20").Value = 76
$ws.Range("F26").Value = 3125
$ws.Range("F27").Value = 225
$ws.Range("F28").Value = 35
$ws.Range("F29").Value = 35
$ws.Range("F30").Value = 433
$ws.Range("F33").Value = 286
$ws.Range("F35").Value = 254

# Sheet 4 (全部类型): rows 36-49 content re-shuffle.
# A new event ("wio夏时之鸢代号鸢Only") was inserted at row 36, pushing the
# subsequent events down by one row (their B/C/D/E/F/G/H/I move to row N+1,
# while column A - a static running index - is left untouched). The old
# "LoveLive！电视动画播放十周年纪念巡演" (sold out) entry that this pushed off the
# end (old row 48) is dropped entirely; rows 49-50 keep their original identity.
$ws = $wb.Worksheets.Item(4)
# Row 36
$ws.Range("B36").Value = '''2024-08-04'
$ws.Range("C36").Value = '广州·wio夏时之鸢代号鸢Only'
$ws.Range("D36").Value = '黄边三横路一街1号 设计殿堂'
$ws.Range("E36").Value = '2024.08.04 10:00-08.04 17:30'
$ws.Range("F36").Value = 226
$ws.Range("G36").Value = 68.8
$ws.Range("H36").Value = 'https://show.bilibili.com/platform/detail.html?id=87434'
$ws.Range("I36").Value = '//i0.hdslb.com/bfs/openplatform/202406/orVoRqXY1718293009879.png'

# Row 37
$ws.Range("B37").Value = '''2024-08-04'
$ws.Range("C37").Value = '广州·星之光动漫嘉年华'
$ws.Range("D37").Value = '钟村镇105国道西侧 广州雄峰城展览中心'
$ws.Range("E37").Value = '2024.08.04 10:00-08.04 17:00'
$ws.Range("F37").Value = 330
$ws.Range("G37").Value = 60
$ws.Range("H37").Value = 'https://show.bilibili.com/platform/detail.html?id=87077'
$ws.Range("I37").Value = '//i2.hdslb.com/bfs/openplatform/202406/hOZ6VVFx1717571239392.jpeg'

# Row 38
$ws.Range("B38").Value = '''2024-08-04'
$ws.Range("C38").Value = '广州·格斗游戏FTGonly'
$ws.Range("D38").Value = '芳村大道下市直街1号信义会馆21栋(近白鹅潭风情酒吧街) 信义会馆-21栋'
$ws.Range("E38").Value = '2024.08.04 10:00-08.04 19:00'
$ws.Range("F38").Value = 18
$ws.Range("G38").Value = 68
$ws.Range("H38").Value = 'https://show.bilibili.com/platform/detail.html?id=87090'
$ws.Range("I38").Value = '//i1.hdslb.com/bfs/openplatform/202406/Vk8sR8Oj1717582522018.png'

# Row 39
$ws.Range("B39").Value = '''2024-08-10'
$ws.Range("C39").Value = '广州·火影忍者only'
$ws.Range("D39").Value = '奥体南路12号 优托邦(奥体旗舰店)'
$ws.Range("E39").Value = '2024.08.10 10:00-08.10 17:00'
$ws.Range("F39").Value = 639
$ws.Range("G39").Value = 60
$ws.Range("H39").Value = 'https://show.bilibili.com/platform/detail.html?id=85704'
$ws.Range("I39").Value = '//i2.hdslb.com/bfs/openplatform/202405/lKOROXve1715763433389.jpeg'

# Row 40
$ws.Range("B40").Value = '''2024-08-10'
$ws.Range("C40").Value = '广州·系统任务：重生之我是音乐一体机！王子健2024巡回演出'
$ws.Range("D40").Value = '广州天河区花城大道89号美食街北二门 SD Livehouse'
$ws.Range("E40").Value = '2024.08.10 20:00-08.10 22:00'
$ws.Range("F40").Value = 44
$ws.Range("G40").Value = 128
$ws.Range("H40").Value = 'https://show.bilibili.com/platform/detail.html?id=87585'
$ws.Range("I40").Value = '//i0.hdslb.com/bfs/openplatform/202406/zIb7ZnHb1718675848837.jpeg'

# Row 41
$ws.Range("B41").Value = '''2024-08-11'
$ws.Range("C41").Value = '广州·咒术回战ONLY'
$ws.Range("D41").Value = '西环路1号 广州岭南会展中心'
$ws.Range("E41").Value = '2024.08.11 10:00-08.11 17:00'
$ws.Range("F41").Value = 174
$ws.Range("G41").Value = 60
$ws.Range("H41").Value = 'https://show.bilibili.com/platform/detail.html?id=87433'
$ws.Range("I41").Value = '//i1.hdslb.com/bfs/openplatform/202406/kNv9yqGn1718350051848.jpeg'

# Row 42
$ws.Range("B42").Value = '''2024-08-11'
$ws.Range("C42").Value = '广州·妖都原神&崩铁only-清凉大作战-'
$ws.Range("D42").Value = '黄边三横路一街1号 设计殿堂'
$ws.Range("E42").Value = '2024.08.11 10:00-08.11 16:30'
$ws.Range("F42").Value = 77
$ws.Range("G42").Value = 60
$ws.Range("H42").Value = 'https://show.bilibili.com/platform/detail.html?id=87321'
$ws.Range("I42").Value = '//i1.hdslb.com/bfs/openplatform/202406/7k54Bi4X1718025336899.jpeg'

# Row 43
$ws.Range("B43").Value = '''2024-08-14'
$ws.Range("C43").Value = '广州·Marcin Patrzalek 2024 《原声之龙》指弹吉他音乐会'
$ws.Range("D43").Value = '海珠同创汇东一街11号（上冲南约11-2） 声音共和Livehouse'
$ws.Range("E43").Value = '2024.08.14 20:00-08.14 21:30'
$ws.Range("F43").Value = 146
$ws.Range("G43").Value = 380
$ws.Range("H43").Value = 'https://show.bilibili.com/platform/detail.html?id=86291'
$ws.Range("I43").Value = '//i1.hdslb.com/bfs/openplatform/202405/vsOXym1L1716546835148.jpeg'

# Row 44
$ws.Range("B44").Value = '''2024-08-16'
$ws.Range("C44").Value = '广州·《最后的莫西干人》-印第安音乐家亚历桑德罗&丛林回响乐队巡演'
$ws.Range("D44").Value = '东风中路299号 广州中山纪念堂'
$ws.Range("E44").Value = '2024.08.16 20:00-08.16 21:30'
$ws.Range("F44").Value = 2
$ws.Range("G44").Value = 380
$ws.Range("H44").Value = 'https://show.bilibili.com/platform/detail.html?id=86143'
$ws.Range("I44").Value = '//i1.hdslb.com/bfs/openplatform/202405/4oOXA1j01716175554059.jpeg'

# Row 45
$ws.Range("B45").Value = '''2024-08-17'
$ws.Range("C45").Value = '广州·COC星火次元云漫创作交流展'
$ws.Range("D45").Value = '黄边三横路一街1号 设计殿堂'
$ws.Range("E45").Value = '2024.08.17 10:00-08.18 17:00'
$ws.Range("F45").Value = 50
$ws.Range("G45").Value = 68
$ws.Range("H45").Value = 'https://show.bilibili.com/platform/detail.html?id=87777'
$ws.Range("I45").Value = '//i1.hdslb.com/bfs/openplatform/202406/nVPxhUFQ1718936306088.jpeg'

# Row 46
$ws.Range("B46").Value = '''2024-08-17'
$ws.Range("C46").Value = '广州·鸟山明作品《龙珠》40周年only纪念展'
$ws.Range("D46").Value = '逸景路462号珠江国际纺织城d区6层 珠江时尚馆'
$ws.Range("E46").Value = '2024.08.17 10:00-08.17 17:30'
$ws.Range("F46").Value = 31
$ws.Range("G46").Value = 88
$ws.Range("H46").Value = 'https://show.bilibili.com/platform/detail.html?id=86780'
$ws.Range("I46").Value = '//i1.hdslb.com/bfs/openplatform/202405/4k7Thger1717147185584.jpeg'

# Row 47
$ws.Range("B47").Value = '''2024-08-18'
$ws.Range("C47").Value = '广州·原神×崩坏×绝区零only'
$ws.Range("D47").Value = '西环路1号 广州岭南会展中心'
$ws.Range("E47").Value = '2024.08.18 10:00-08.18 17:00'
$ws.Range("F47").Value = 568
$ws.Range("G47").Value = 60
$ws.Range("H47").Value = 'https://show.bilibili.com/platform/detail.html?id=87025'
$ws.Range("I47").Value = '//i0.hdslb.com/bfs/openplatform/202405/lsOq4H701717169339283.png'

# Row 48
$ws.Range("B48").Value = '''2024-08-23'
$ws.Range("C48").Value = '广州·LoveLiveOnly'
$ws.Range("D48").Value = '芳村大道下市直街1号信义会馆21栋(近白鹅潭风情酒吧街) 信义会馆-21栋'
$ws.Range("E48").Value = '2024.08.23 10:00-08.23 19:00'
$ws.Range("F48").Value = 225
$ws.Range("G48").Value = 68.8
$ws.Range("H48").Value = 'https://show.bilibili.com/platform/detail.html?id=87033'
$ws.Range("I48").Value = '//i2.hdslb.com/bfs/openplatform/202406/a8shiH411717579829497.jpeg'

# Row 49
$ws.Range("B49").Value = '''2024-08-27'
$ws.Range("C49").Value = '广州·25时主题同人茶会×晓山瑞希生日会'
$ws.Range("D49").Value = '黄边地铁B出口黄边美食广场1层 胡桃里音乐馆(黄边店)'
$ws.Range("E49").Value = '2024.08.27 10:00-08.27 16:30'
$ws.Range("F49").Value = 28
$ws.Range("G49").Value = 58
$ws.Range("H49").Value = 'https://show.bilibili.com/platform/detail.html?id=87815'
$ws.Range("I49").Value = '//i1.hdslb.com/bfs/openplatform/202406/rzS5X2Ko1718735908971.png'

